$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new PR reference value + hyperlink in B17
$url = "https://github.com/dhavalkeerthi/MRIInterns2026A/pull/5"
$ws.Range("B17").Value = $url
$ws.Hyperlinks.Add($ws.Range("B17"), $url)

# Match the style used by the existing hyperlink cell (B20)
$ws.Range("B17").Style = $ws.Range("B20").Style

# Update the view state to match: scrolled so row 13 is the top row,
# and the active selection on D17
$excel.ActiveWindow.ScrollRow = 13
$ws.Range("D17").Select()
